# Fruta / hortaliza, semanal
# Insert a new weekly record row right above the old row 397 (pushing all
# following rows down by one), duplicating the data that used to live in
# row 397 but with the new week's date (Fecha) and Volumen values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 397; everything below shifts down.
$ws.Rows.Item(397).Insert()

# Populate the newly inserted row 397 with the same record that used to be
# there (now at row 398), except for the updated Fecha (D) and Volumen (J).
$ws.Range("A397").Value = 11
$ws.Range("B397").Value = "Vega Monumental Concepción"
$ws.Range("C397").Value = "Bíobío"
$ws.Range("D397").Value = 44931
$ws.Range("E397").Value = 8
$ws.Range("F397").Value = 100114014
$ws.Range("G397").Value = "Betarraga"
$ws.Range("H397").Value = "Sin especificar"
$ws.Range("I397").Value = "Primera"
$ws.Range("J397").Value = 1150
$ws.Range("K397").Value = 600
$ws.Range("L397").Value = 650
$ws.Range("M397").Value = 624
$ws.Range("N397").Value = "$/paquete 5 unidades"
$ws.Range("O397").Value = "Región Metropolitana"
$ws.Range("P397").Value = 125
$ws.Range("Q397").Value = 5
$ws.Range("R397").Value = "Hortaliza"

# Make sure the Fecha cell keeps the date style used throughout column D.
$ws.Range("D397").NumberFormat = $ws.Range("D398").NumberFormat
